$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F2:F13 updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 79
$ws1.Range("F3").Value = 806
$ws1.Range("F4").Value = 44
$ws1.Range("F5").Value = 62
$ws1.Range("F6").Value = 111
$ws1.Range("F7").Value = 342
$ws1.Range("F8").Value = 4420
$ws1.Range("F9").Value = 95
$ws1.Range("F10").Value = 4973
$ws1.Range("F11").Value = 558
$ws1.Range("F12").Value = 1251
$ws1.Range("F13").Value = 86

# Sheet "全部类型" (sheet4): F2:F14 updates (F7 unchanged - belongs to "演出" entry)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 79
$ws4.Range("F3").Value = 806
$ws4.Range("F4").Value = 44
$ws4.Range("F5").Value = 62
$ws4.Range("F6").Value = 111
$ws4.Range("F8").Value = 342
$ws4.Range("F9").Value = 4420
$ws4.Range("F10").Value = 95
$ws4.Range("F11").Value = 4973
$ws4.Range("F12").Value = 558
$ws4.Range("F13").Value = 1251
$ws4.Range("F14").Value = 86
